# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff" everywhere it
#   appears (Overview zh-cn/de-de status columns + per-language Status column).
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to the new handoff-report generation time.
# - Widen the status/datetime columns so the longer "Ready for handoff" text
#   (and the datetime stamps) aren't truncated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps refreshed for the new handoff generation ---
$overview.Range("G2").Value = "2016-08-12 05:01:45"
$zhcn.Range("H2").Value = "2016-08-12 05:01:39"
$dede.Range("H2").Value = "2016-08-12 05:01:45"

# --- Widen the status/datetime columns to fit the new content ---
$overview.Columns.Item(5).ColumnWidth = 16.38265482584637
$overview.Columns.Item(6).ColumnWidth = 16.38265482584637
$zhcn.Columns.Item(3).ColumnWidth = 16.38265482584637
$dede.Columns.Item(3).ColumnWidth = 16.38265482584637
